$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.655.77"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.707.79"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "673.44"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.03"
$ws.Range("E6").Value = "  +2.52%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.498"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.13"
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.88"
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").Value = "3.704.10"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Value = "69.673.36"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("E17").Value = "  +2.82%  "
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "473.40"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "3.857.86"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("E24").Value = "  +5.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.02"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.17"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.167"
$ws.Range("E32").Value = "  +6.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").Value = "3.697.14"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("E36").Value = "  +5.03%  "
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.27"
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0914"
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "174.18"
$ws.Range("E42").Value = "  +4.80%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.03"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("E45").Value = "  +2.51%  "
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("E48").Value = "  +3.74%  "
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.88"
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("E51").Value = "  +1.74%  "
